# Fruta / hortaliza, semanal
# Shifts the weekly price records in rows 198-310 down by one row (columns
# D, J, K, L, M, P), inserts a brand-new latest record in row 198, and
# appends the record that used to be in row 310 as the new row 311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 198
$lastRow  = 310

# --- 1. Snapshot the columns that move (D,J,K,L,M,P) for every row in the
#        block, before we start overwriting anything. -----------------------
$snapD = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value()
    $snapJ[$r] = $ws.Cells.Item($r, 10).Value()
    $snapK[$r] = $ws.Cells.Item($r, 11).Value()
    $snapL[$r] = $ws.Cells.Item($r, 12).Value()
    $snapM[$r] = $ws.Cells.Item($r, 13).Value()
    $snapP[$r] = $ws.Cells.Item($r, 16).Value()
}

# --- 2. Snapshot the whole of the last row (310) -- it becomes row 311 -----
$lastRowVals = @{}
for ($c = 1; $c -le 18; $c++) {
    $lastRowVals[$c] = $ws.Cells.Item($lastRow, $c).Value()
}

# --- 3. Append the new row 311 = old row 310, verbatim ---------------------
$newRow = $lastRow + 1
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newRow, $c).Value = $lastRowVals[$c]
}
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 4. Push every row's D/J/K/L/M/P down into the next row ----------------
#        (row 310 gets what used to be in row 309, ... row 199 gets what
#        used to be in row 198).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value  = $snapD[$src]
    $ws.Cells.Item($r, 10).Value = $snapJ[$src]
    $ws.Cells.Item($r, 11).Value = $snapK[$src]
    $ws.Cells.Item($r, 12).Value = $snapL[$src]
    $ws.Cells.Item($r, 13).Value = $snapM[$src]
    $ws.Cells.Item($r, 16).Value = $snapP[$src]
}

# --- 5. Row 198 becomes a brand-new record (date + volume change only;
#        min/max/avg price and Precio $/Kg are unchanged from before). -----
$ws.Cells.Item($firstRow, 4).Value  = 44806
$ws.Cells.Item($firstRow, 10).Value = 2000
$ws.Cells.Item($firstRow, 11).Value = $snapK[$firstRow]
$ws.Cells.Item($firstRow, 12).Value = $snapL[$firstRow]
$ws.Cells.Item($firstRow, 13).Value = $snapM[$firstRow]
$ws.Cells.Item($firstRow, 16).Value = $snapP[$firstRow]

Write-Output "done"
